$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New transaction rows (Passing no meja / relayout Integrasi)
$ws.Range("A43").Value = 43283.715720532404
$ws.Range("B43").Value = "#ID02034"
$ws.Range("C43").Value = "Kacamata"
$ws.Range("D43").Value = 142780.0

$ws.Range("A44").Value = 43284.37010010417
$ws.Range("B44").Value = "#ID03035"
$ws.Range("C44").Value = "Oto Bento"
$ws.Range("D44").Value = 27500.000000000004

$ws.Range("A45").Value = 43284.372528912034
$ws.Range("B45").Value = "#ID03036"
$ws.Range("C45").Value = "PHD"
$ws.Range("D45").Value = 167200.0
